$wb = $excel.ActiveWorkbook

$oldUuid = "d0ca1dd6-2ade-4e81-a460-1e79887cd4d5"
$newUuid = "fa7653dc-e02a-4135-88b8-2786525df508"
$oldHash = "4f97ae72bdc331ac2cf7ad237d892bf152f05406"
$newHash = "1221e1959d851dd760da248d07bae31863266a27"

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newUuid.md"
$wsOverview.Range("G2").Value = "2016-08-22 17:02:22"

# Update the B2 hyperlink display text while keeping its original target.
$wsOverview.Range("A1:Z100").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49dfb872f10c5b0168fd6213829b3300212f820a/e2e/$oldUuid.md",
    "",
    "",
    "e2e\$newUuid.md"
) | Out-Null

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newUuid.md"
$wsZhCn.Range("G2").Value = "$newUuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-22 17:02:17"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# Remove the "Latest Target File" (I2) hyperlink while keeping the A2 one,
# refreshed with the new display text but the original target address.
$wsZhCn.Range("A1:Z100").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49dfb872f10c5b0168fd6213829b3300212f820a/e2e/$oldUuid.md",
    "",
    "",
    "$newUuid.md"
) | Out-Null

$wsZhCn.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZhCn.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newUuid.md"
$wsDeDe.Range("G2").Value = "$newUuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-22 17:02:22"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Range("A1:Z100").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49dfb872f10c5b0168fd6213829b3300212f820a/e2e/$oldUuid.md",
    "",
    "",
    "$newUuid.md"
) | Out-Null

$wsDeDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDeDe.Columns.Item(10).ColumnWidth = 21.7054770333426
